$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 453 (shifts existing rows 453:482 down to 454:483)
$ws.Rows.Item(453).Insert()

# Populate the newly inserted row 453 with the new data record
$ws.Range("A453").Value = 3
$ws.Range("B453").Value = "Femacal de La Calera"
$ws.Range("C453").Value = "Coquimbo"
$ws.Range("D453").Value = 44610
$ws.Range("E453").Value = 5
$ws.Range("F453").Value = 100114001
$ws.Range("G453").Value = "Papa"
$ws.Range("H453").Value = "Rosara"
$ws.Range("I453").Value = "1a (cosecha)"
$ws.Range("J453").Value = 340
$ws.Range("K453").Value = 7000
$ws.Range("L453").Value = 7500
$ws.Range("M453").Value = 7235
$ws.Range("N453").Value = "`$/saco 25 kilos"
$ws.Range("O453").Value = "Provincia de Quillota"
$ws.Range("P453").Value = 289
$ws.Range("Q453").Value = 25
$ws.Range("R453").Value = "Hortaliza"
